$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Hunk 1: merge the two runs "m" + "öchte ich" (the occurrence that follows
# "Besucher der Lernpattform") into a single run "möchte ich".
# There is an earlier, untouched "m" / "öchte ich" pair (template row with
# "<Rolle>") further up in the document, so we anchor the search right after
# the "Besucher der Lernpattform" text to land on the correct occurrence.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Besucher der Lernpattform")
$rng.Collapse(0)
$rng.Find.Execute("möchte ich")
# Assigning the exact same text the range already has is a no-op for the
# underlying run layout, so bounce through a distinct placeholder string
# first to force the engine to actually collapse "m" + "öchte ich" into one
# run before writing the final text back.
$rng.Text = "möchte ichXX"
$rng.Text = "möchte ich"

# ---------------------------------------------------------------------------
# Hunk 2: "Kurse und Aufgaben erstellen können" -> "Kurse und Aufgaben erstellen"
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Kurse und Aufgaben erstellen können")
$rng.Text = "Kurse und Aufgaben erstellen"

# ---------------------------------------------------------------------------
# Hunk 3: "diese meinen Kursteilnehmern anbieten zu können" is split into two
# runs: "diese meinen Kursteilnehmern " (preserve trailing space) + "anzubieten"
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("diese meinen Kursteilnehmern anbieten zu können")
$rng.Text = "diese meinen Kursteilnehmern "
$rngA = $d.Range($rng.Start, $rng.End)
$rng.Collapse(0)
$rng.InsertAfter("anzubieten")
$rngB = $d.Range($rng.Start, $rng.End)
# Force the two freshly written runs to stay distinct instead of being
# coalesced back together by the engine's adjacent-identical-run merge.
$rngA.Bold = 1
$rngA.Bold = 0
$rngB.Bold = 1
$rngB.Bold = 0

# ---------------------------------------------------------------------------
# Hunk 4: "mich für Kurse anmelden können" -> "mich für Kurse anmelden"
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("mich für Kurse anmelden können")
$rng.Text = "mich für Kurse anmelden"

# ---------------------------------------------------------------------------
# Hunk 5: "nach der Anmeldung meine Kurse sehen können und neue Kurse
# hinzufügen können" is split into three runs: "nach der Anmeldung meine
# Kurse sehen" + " " (preserve) + "und neue Kurse hinzufügen"
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("nach der Anmeldung meine Kurse sehen können und neue Kurse hinzufügen können")
$rng.Text = "nach der Anmeldung meine Kurse sehen"
$rngA = $d.Range($rng.Start, $rng.End)
$rng.Collapse(0)
$rng.InsertAfter(" ")
$rngB = $d.Range($rng.Start, $rng.End)
$rng.Collapse(0)
$rng.InsertAfter("und neue Kurse hinzufügen")
$rngC = $d.Range($rng.Start, $rng.End)
$rngA.Bold = 1
$rngA.Bold = 0
$rngB.Bold = 1
$rngB.Bold = 0
$rngC.Bold = 1
$rngC.Bold = 0

# ---------------------------------------------------------------------------
# Hunk 6: "die Rolle eines Kursteilnehmers einnehmen können" is split into
# three runs: "die" + " auch" (preserve) + " Rolle eines Kursteilnehmers
# einnehmen" (preserve)
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("die Rolle eines Kursteilnehmers einnehmen können")
$rng.Text = "die"
$rngA = $d.Range($rng.Start, $rng.End)
$rng.Collapse(0)
$rng.InsertAfter(" auch")
$rngB = $d.Range($rng.Start, $rng.End)
$rng.Collapse(0)
$rng.InsertAfter(" Rolle eines Kursteilnehmers einnehmen")
$rngC = $d.Range($rng.Start, $rng.End)
$rngA.Bold = 1
$rngA.Bold = 0
$rngB.Bold = 1
$rngB.Bold = 0
$rngC.Bold = 1
$rngC.Bold = 0
